$wb = $excel.ActiveWorkbook

# --- Sheet: Rocket Parameters (Mass) ---
$ws = $wb.Worksheets.Item("Rocket Parameters (Mass)")

# B6: 1 -> 2
$ws.Range("B6").Value = 2

# B8: was formula =P2 -> now a literal value 2.5
$ws.Range("B8").Value = 2.5

# B13: was formula =P5 -> now a literal value 1.5
$ws.Range("B13").Value = 1.5

# New row 15: Fudge Factor
$ws.Range("A15").Value = "Fudge"
$ws.Range("B15").Formula = "=SUM(B3:B14)*E15"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = "kg"
$ws.Range("D15").Value = "Factor"
$ws.Range("E15").Value = 0.15
$ws.Range("E15").NumberFormat = "0%"

# Update selection to B16
$ws.Range("B16").Select() | Out-Null

# --- Sheet: Engine Parameters ---
$ws3 = $wb.Worksheets.Item("Engine Parameters")
$ws3.Range("C6").Value = 4.99
$ws3.Range("C7").Select() | Out-Null

# --- Sheet: Propellant Parameters (Tanks) ---
$ws4 = $wb.Worksheets.Item("Propellant Parameters (Tanks)")
$ws4.Range("J7").Value = 1600
$ws4.Range("J9").Value = 920
$ws4.Range("N20").Select() | Out-Null
